$wb = $excel.ActiveWorkbook

# Sheet4 = 景点门票 (attraction tickets)
$ws4 = $wb.Worksheets.Item("景点门票")

# Add "total" row under the existing data
$ws4.Range("A5").Value = "total"
$ws4.Range("C5").Formula = "=SUM(C2:C4)"

# Update selection on sheet4 to A6
$ws4.Range("A6").Select()

# Sheet1 = 行程 (itinerary) - recalc B18 total and move selection to B18
$ws1 = $wb.Worksheets.Item("行程")
$ws1.Range("B18").Select()

$wb.Application.Calculate()
